# Estadisticos Matutinos 15 Oct
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheets "Estadisticos 1P" and "Estadisticos Final": columns D:H, rows 2-7
# ---------------------------------------------------------------------------
$sheet1Names = @("Estadisticos 1P", "Estadisticos Final")
foreach ($name in $sheet1Names) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("D2").Value = 0
    $ws.Range("E2").Value = 9
    $ws.Range("F2").Value = 22
    $ws.Range("G2").Value = 70.97
    $ws.Range("H2").Value = 6.6

    $ws.Range("D3").Value = 0
    $ws.Range("E3").Value = 7
    $ws.Range("F3").Value = 29
    $ws.Range("G3").Value = 80.56
    $ws.Range("H3").Value = 7.5

    $ws.Range("D4").Value = 2
    $ws.Range("E4").Value = 7
    $ws.Range("F4").Value = 34
    $ws.Range("G4").Value = 79.06999999999999
    $ws.Range("H4").Value = 7.5

    $ws.Range("D5").Value = 0
    $ws.Range("E5").Value = 6
    $ws.Range("F5").Value = 38
    $ws.Range("G5").Value = 86.36
    $ws.Range("H5").Value = 7

    $ws.Range("D6").Value = 0
    $ws.Range("E6").Value = 3
    $ws.Range("F6").Value = 21
    $ws.Range("G6").Value = 87.5
    $ws.Range("H6").Value = 7.1

    $ws.Range("D7").Value = 0
    $ws.Range("E7").Value = 3
    $ws.Range("F7").Value = 23
    $ws.Range("G7").Value = 88.45999999999999
    $ws.Range("H7").Value = 6.7
}

# ---------------------------------------------------------------------------
# Sheet "Estadisticos 2P": only column E, rows 2-7, changes
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("E2").Value = 31
$ws2.Range("E3").Value = 36
$ws2.Range("E4").Value = 41
$ws2.Range("E5").Value = 44
$ws2.Range("E6").Value = 24
$ws2.Range("E7").Value = 26

# ---------------------------------------------------------------------------
# Sheet "Rescatables": new rows 2-5 with rescue-exam candidate data
# Strings are written column-by-column (B, then C, then D) to reproduce the
# shared-string insertion order of the authored workbook.
# ---------------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

$ws4.Range("A2").Value = 21330051920201
$ws4.Range("A3").Value = 21330051920288
$ws4.Range("A4").Value = 21330051920327
$ws4.Range("A5").Value = 21330051920243

$ws4.Range("B2").Value = "OFICIAL"
$ws4.Range("B3").Value = "VELAZQUEZ"
$ws4.Range("B4").Value = "TEXCAHUA"
$ws4.Range("B5").Value = "PERALTA"

$ws4.Range("C2").Value = "TZOMPAXTLE"
$ws4.Range("C3").Value = "TEXCAHUA"
$ws4.Range("C4").Value = "MARTINEZ"
$ws4.Range("C5").Value = "HERNANDEZ"

$ws4.Range("D2").Value = "CLEMENTE"
$ws4.Range("D3").Value = "NANCY PAOLA"
$ws4.Range("D4").Value = "MARISOL"
$ws4.Range("D5").Value = "MICHELLE GUADALUPE"

$ws4.Range("E2").Value = "INGLÉS I"
$ws4.Range("E3").Value = "INGLÉS I"
$ws4.Range("E4").Value = "INGLÉS I"
$ws4.Range("E5").Value = "INGLÉS I"

$ws4.Range("F2").Value = "1BM"
$ws4.Range("F3").Value = "1DM"
$ws4.Range("F4").Value = "1EM"
$ws4.Range("F5").Value = "1CM"

$ws4.Range("G2").Value = 6
$ws4.Range("G3").Value = 6
$ws4.Range("G4").Value = 6
$ws4.Range("G5").Value = 6
